$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-17 (Fixture, Pick, Competition, Time, Win Info, Confidence %, Odds)
$data = @(
    @("Chelsea vs Sunderland", "Chelsea", "England", "25th Oct 15:00", "82 of 106 Tips", "77", "1.45"),
    @("Brentford vs Liverpool", "Liverpool", "England", "25th Oct 20:00", "63 of 88 Tips", "72", "1.83"),
    @("Newcastle vs Fulham", "Newcastle", "England", "25th Oct 15:00", "57 of 84 Tips", "68", "1.61"),
    @("Borussia Monchengladbach vs Bayern Munich", "Bayern Munich", "Germany", "25th Oct 14:45", "50 of 51 Tips", "98", "1.28"),
    @("Brest vs PSG", "PSG", "France", "25th Oct 16:00", "42 of 45 Tips", "93", "1.35"),
    @("Athletic Bilbao vs Getafe", "Athletic Bilbao", "Spain", "25th Oct 17:30", "39 of 44 Tips", "89", "1.67"),
    @("Borussia Dortmund vs FC Cologne", "Borussia Dortmund", "Germany", "25th Oct 17:30", "36 of 36 Tips", "100", "1.42"),
    @("Monaco vs Toulouse", "Monaco", "France", "25th Oct 18:00", "33 of 37 Tips", "89", "1.92"),
    @("Espanyol vs Elche", "Espanyol", "Spain", "25th Oct 15:15", "30 of 36 Tips", "83", "2.05"),
    @("Arsenal vs Crystal Palace", "Arsenal", "England", "26th Oct 14:00", "23 of 28 Tips", "82", "1.40"),
    @("Valencia vs Villarreal", "Villarreal", "Spain", "25th Oct 20:00", "19 of 29 Tips", "66", "2.05"),
    @("Bournemouth vs Nottingham Forest", "Bournemouth", "England", "26th Oct 14:00", "19 of 23 Tips", "83", "1.83"),
    @("Aston Villa vs Man City", "Man City", "England", "26th Oct 14:00", "18 of 24 Tips", "75", "1.80"),
    @("Cremonese vs Atalanta", "Atalanta", "Italy", "25th Oct 19:45", "17 of 23 Tips", "74", "1.57"),
    @("Real Madrid vs Barcelona", "Real Madrid", "Spain", "26th Oct 15:15", "16 of 24 Tips", "67", "2.05"),
    @("Fluminense vs Internacional", "Fluminense", "Brazil", "25th Oct 21:30", "16 of 17 Tips", "94", "1.83")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
    $ws.Cells.Item($row, 7).Value = $rowData[6]
}

# Remove old rows 18 and 19 (formerly Valencia/Villarreal and Bournemouth rows)
$ws.Range("A18:G19").Delete() | Out-Null

# Clear the old H20 formula cell and set the new H18 formula
$ws.Range("H20").ClearContents() | Out-Null
$ws.Range("H18").Formula = "=AVERAGE(F2:F17)"
